$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Problem"
$ws.Cells.Item(1, 2).Value = "Description"
$ws.Cells.Item(1, 3).Value = "Date_Created"
$ws.Cells.Item(1, 4).Value = "Start Date"
$ws.Cells.Item(1, 5).Value = "Email attached file"

$ws.Cells.Item(2, 1).Value = "ללכת לאכול בפארק"
$ws.Cells.Item(2, 2).Value = "רעב רצחחחחחחחח חיייייייבבבב אוכלללללל"
$ws.Cells.Item(2, 3).Value = "2021-01-20 20:35:26.623884+00:00"
$ws.Cells.Item(2, 4).Value = "2021-01-21 22:35:00+00:00"
$ws.Cells.Item(2, 5).Value = "Emails_Files_Tasks/Ht.png"

$ws.Cells.Item(3, 1).Value = "סו ושמחו , בשמחת תורה"
$ws.Cells.Item(3, 2).Value = "הייייייי אי ניייד סומבאדייי היייי"
$ws.Cells.Item(3, 3).Value = "2021-01-20 20:46:39.909970+00:00"
$ws.Cells.Item(3, 4).Value = "2021-01-20 22:46:00+00:00"
$ws.Cells.Item(3, 5).Value = "Emails_Files_Tasks/Ht_XpPHx39.png"

$ws.Cells.Item(4, 1).Value = "Tomer Test"
$ws.Cells.Item(4, 2).Value = "גדשגדשגדשגדשגדגדשגדש"
$ws.Cells.Item(4, 3).Value = "2021-01-20 20:47:15.406410+00:00"
$ws.Cells.Item(4, 4).Value = "2021-01-09 22:47:00+00:00"
$ws.Cells.Item(4, 5).Value = "Emails_Files_Tasks/ClearBitAPI_VwWZRui.PNG"

$ws.Cells.Item(5, 1).Value = "Hello World"
$ws.Cells.Item(5, 2).Value = "hello hello hellp"
$ws.Cells.Item(5, 3).Value = "2021-01-26 18:34:36.136413+00:00"
$ws.Cells.Item(5, 4).Value = "2021-02-02 20:33:00+00:00"
$ws.Cells.Item(5, 5).Value = "Emails_Files_Tasks/ClearBitAPI_W55aLN4.PNG"

$ws.Cells.Item(6, 1).Value = "DUDU topaz"
$ws.Cells.Item(6, 2).Value = "wassuppppp"
$ws.Cells.Item(6, 3).Value = "2021-01-27 23:16:57.643539+00:00"
$ws.Cells.Item(6, 4).Value = "2021-01-01 01:16:00+00:00"
$ws.Cells.Item(6, 5).Value = "No File Attached"

$ws.Cells.Item(7, 1).Value = "eat alot of food"
$ws.Cells.Item(7, 2).Value = "PIZZAA PIZAAA PIZZZA"
$ws.Cells.Item(7, 3).Value = "2021-01-27 23:21:29.633551+00:00"
$ws.Cells.Item(7, 4).Value = "2021-01-02 01:21:00+00:00"
$ws.Cells.Item(7, 5).Value = "No File Attached"

$ws.Cells.Item(8, 1).Value = "Eat Eat"
$ws.Cells.Item(8, 2).Value = "Eattttttttttttttttttt"
$ws.Cells.Item(8, 3).Value = "2021-01-28 13:41:34.927048+00:00"
$ws.Cells.Item(8, 4).Value = "2021-01-04 15:41:00+00:00"
$ws.Cells.Item(8, 5).Value = "No File Attached"

$ws.Cells.Item(9, 1).Value = "לאכולללללללללללללללללללל"
$ws.Cells.Item(9, 2).Value = "לאכול לאכול מלא אוכל לאכול"
$ws.Cells.Item(9, 3).Value = "2021-01-28 13:46:37.974221+00:00"
$ws.Cells.Item(9, 4).Value = "2021-01-11 15:44:00+00:00"
$ws.Cells.Item(9, 5).Value = "No File Attached"

$ws.Range("A1:E1").Font.Bold = $true
